# "data modified in IPL folder"
# - Insert two new columns "ownTeam"/"oppTeam" after column C (shifting the
#   old batsman..sr columns from D:I to F:K).
# - Re-order the match rows chronologically and fill in the new team columns
#   (the sheet always shows Royal Challengers Bangalore as ownTeam, the
#   opponent for that match as oppTeam) for every row.
# All values on this sheet are plain text (original file used t="str" for
# every cell, including numeric-looking ones like run counts / strike
# rates), so we force Text number format before writing so Excel does not
# silently reinterpret values such as "4" or "200.00" as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (ownTeam, oppTeam) at D and E, shifting old D:I to F:K
$ws.Columns("D:E").Insert()

# Force the whole used range to Text format so values are stored as strings, not numbers
$ws.Range("A1:K6").NumberFormat = "@"

$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 28 2020"
$ws.Range("C2").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Chris Morris "
$ws.Range("G2").Value = "4"
$ws.Range("H2").Value = "2"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "200.00"

$ws.Range("A3").Value = " Sharjah"
$ws.Range("B3").Value = " October 15 2020"
$ws.Range("C3").Value = "Kings XI won by 8 wickets"
$ws.Range("D3").Value = "Royal Challengers Bangalore"
$ws.Range("E3").Value = "Kings XI Punjab"
$ws.Range("F3").Value = "Chris Morris "
$ws.Range("G3").Value = "25"
$ws.Range("H3").Value = "8"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "3"
$ws.Range("K3").Value = "312.50"

$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 31 2020"
$ws.Range("C4").Value = "Sunrisers won by 5 wickets (with 35 balls remaining)"
$ws.Range("D4").Value = "Royal Challengers Bangalore"
$ws.Range("E4").Value = "Sunrisers Hyderabad"
$ws.Range("F4").Value = "Chris Morris "
$ws.Range("G4").Value = "3"
$ws.Range("H4").Value = "4"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "75.00"

$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " November 02 2020"
$ws.Range("C5").Value = "Capitals won by 6 wickets (with 6 balls remaining)"
$ws.Range("D5").Value = "Royal Challengers Bangalore"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Chris Morris "
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "2"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "0.00"

$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 25 2020"
$ws.Range("C6").Value = "Super Kings won by 8 wickets (with 8 balls remaining)"
$ws.Range("D6").Value = "Royal Challengers Bangalore"
$ws.Range("E6").Value = "Chennai Super Kings"
$ws.Range("F6").Value = "Chris Morris "
$ws.Range("G6").Value = "2"
$ws.Range("H6").Value = "5"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "40.00"
